$d = $word.ActiveDocument

# The three logo pictures living in the headers/footers need their
# picture "name" swapped (image1.png <-> image2.png for the Pearson
# logo, image2.jpg -> image1.jpg for the BTEC logo). We identify each
# picture by its (stable) AlternativeText/description rather than by
# a hard-coded header/footer index, so the script is robust regardless
# of how many sections the document has.

function Update-LogoName {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $desc = $shp.AlternativeText

        if ($desc -eq "BTec_Logo-Orange") {
            if ($shp.Name -ne "image1.jpg") {
                $shp.Name = "image1.jpg"
            }
        }
        elseif ($desc -like "*PearsonLogo.png") {
            if ($shp.Name -ne "image2.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($h = 1; $h -le $sec.Headers.Count; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            Update-LogoName $hdr.Range.InlineShapes
        }
    }

    for ($f = 1; $f -le $sec.Footers.Count; $f++) {
        $ftr = $sec.Footers.Item($f)
        if ($ftr.Exists) {
            Update-LogoName $ftr.Range.InlineShapes
        }
    }
}
